$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

$ws.Range("C1").Value = "gpu"
$ws.Range("C2").Value = 0
$ws.Range("C3").Value = 1
$ws.Range("C4").Value = 0
$ws.Range("C5").Value = 1
$ws.Range("C6").Value = 0
$ws.Range("C7").Value = 1
$ws.Range("C8").Value = 0
$ws.Range("C9").Value = 1

$ws.Range("C10").Select()
